# Apply the "improve MAPE and include on reporting" update to the
# CustomerService MONTHLY_inventory report on sheet "Main".
#
# 1. Clear the (unused) number-format style that had been stamped on every
#    department_name cell in column B (rows 2-77) -- it never carried an
#    actual custom format, so resetting the style to Normal removes the
#    stray cellXfs entry.
# 2. Refresh the inventory_cases_end_of_month (column E) figures that
#    changed with the improved MAPE calculation.
# 3. Scroll the sheet view down so row 52 is the first visible row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the leftover style from column B (department_name) ---------
$ws.Range("B2:B77").Style = "Normal"

# --- 2. Updated inventory figures in column E -----------------------------
$updates = @{
    4  = 74
    8  = 30
    10 = 103
    13 = 50
    28 = 120
    36 = 25
    41 = 118
    42 = 2
    53 = 24
    63 = 122
    70 = 193
    71 = 100
    73 = 39
    74 = 7
    77 = 105
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}

# --- 3. Scroll the active sheet view so A52 is the top-left cell ---------
$excel.ActiveWindow.TopLeftCell = $ws.Range("A52")
